# 自动更新Excel文件 - 2025-11-03 23:12:17
#
# Daily housekeeping pass over the tracking sheet:
#   D = 总天 (total days allotted)
#   E = 剩余 (days remaining)
#   F = 开始时间 (start date, yyyymmdd)
#
# "Today" has advanced one day since the sheet was last refreshed, so every
# live row's remaining-day counter (E) ticks down by 1. A row whose counter
# had already reached 1 (i.e. it would hit 0 today) is instead renewed: its
# remaining days reset back to the full allotment (D) and its start date is
# rolled forward to the new "today" (2025-11-04).
#
# Rows whose start date isn't a clean 8-digit yyyymmdd value are data-entry
# errors and are left untouched, same as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = 20251104
$firstRow = 2
$lastRow = 99

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $startDate = $ws.Cells.Item($row, 6).Value2
    $startDateText = [string]$startDate

    if ($startDateText.Length -ne 8) {
        # Malformed/blank start date (e.g. "202510929") - skip, don't touch.
        continue
    }

    $totalDays = $ws.Cells.Item($row, 4).Value2
    $remaining = $ws.Cells.Item($row, 5).Value2

    if ($remaining -eq 1) {
        # Counter would expire today - renew the entry.
        $ws.Cells.Item($row, 5).Value = $totalDays
        $ws.Cells.Item($row, 6).Value = $today
    } else {
        # Normal daily decrement.
        $ws.Cells.Item($row, 5).Value = $remaining - 1
    }
}
